$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G17 value (hours), which will cascade formula recalculation for
# H17, G36, and H36 automatically.
$ws.Range("G17").Value = 2

# Widen column C slightly (engine quantizes column width to whole pixels,
# so 44 is the closest achievable input to the target stored width).
$ws.Columns("C").ColumnWidth = 44

# Update view state: scroll so row 3 is at top, and select G18.
$excel.ActiveWindow.ScrollRow = 3
[void]$ws.Range("G18").Select()
